$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 5185.7144
$ws.Range("J52").Value = 5950
$ws.Range("L52").Value = 17850
$ws.Range("N52").Value = -18170

$ws.Range("H137").Value = 1810.5172
$ws.Range("I137").Value = 1110.0555
$ws.Range("J137").Value = 2956.7273
$ws.Range("K137").Value = 3330.1665
$ws.Range("L137").Value = 8870.1819
$ws.Range("M137").Value = -780.1664999999998
$ws.Range("N137").Value = -13970.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2877.0176
$ws.Range("I32").Value = 3042.7292
$ws.Range("K32").Value = 3042.7292
$ws.Range("M32").Value = -2755.7292

$ws.Range("H43").Value = 7999.5
$ws.Range("I43").Value = 7999
$ws.Range("K43").Value = 7999
$ws.Range("M43").Value = -7686

$ws.Range("H61").Value = 1012
$ws.Range("I61").Value = 1012
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1012
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -800
$ws.Range("N61").ClearContents()

$ws.Range("H132").Value = 2410.5908
$ws.Range("I132").Value = 2170.9443
$ws.Range("K132").Value = 6512.8329
$ws.Range("M132").Value = -3982.8329

$ws.Range("H136").Value = 1012
$ws.Range("I136").Value = 1012
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3036
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -486
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 41668050
$ws.Range("I99").Value = 41668050
$ws.Range("K99").Value = 41668050
$ws.Range("M99").Value = -41666552

$ws.Range("H105").Value = 83335790
$ws.Range("I105").Value = 125002580
$ws.Range("K105").Value = 125002580
$ws.Range("M105").Value = -125000833

$ws.Range("H132").Value = 36259.668
$ws.Range("J132").Value = 36259.668
$ws.Range("L132").Value = 36259.668
$ws.Range("N132").Value = -46379.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4997535
$ws.Range("I6").Value = 6660046.5
$ws.Range("K6").Value = 6660046.5
$ws.Range("M6").Value = -6659933.5

$ws.Range("H7").Value = 291.41666
$ws.Range("I7").Value = 333.6
$ws.Range("J7").Value = 80.5
$ws.Range("K7").Value = 333.6
$ws.Range("L7").Value = 80.5
$ws.Range("M7").Value = -220.6
$ws.Range("N7").Value = -306.5

$ws.Range("H31").Value = 991.4912
$ws.Range("I31").Value = 679.95
$ws.Range("J31").Value = 1724.5294
$ws.Range("K31").Value = 679.95
$ws.Range("L31").Value = 1724.5294
$ws.Range("M31").Value = -384.95
$ws.Range("N31").Value = -2314.5294

$ws.Range("H34").Value = 991.4912
$ws.Range("I34").Value = 679.95
$ws.Range("J34").Value = 1724.5294
$ws.Range("K34").Value = 679.95
$ws.Range("L34").Value = 1724.5294
$ws.Range("M34").Value = -477.95
$ws.Range("N34").Value = -2128.5294

$ws.Range("H62").Value = 6454516.5
$ws.Range("I62").Value = 2988
$ws.Range("J62").Value = 33335884
$ws.Range("K62").Value = 2988
$ws.Range("L62").Value = 33335884
$ws.Range("M62").Value = -2364
$ws.Range("N62").Value = -33337132

$ws.Range("H65").Value = 6454516.5
$ws.Range("I65").Value = 2988
$ws.Range("J65").Value = 33335884
$ws.Range("K65").Value = 14940
$ws.Range("L65").Value = 166679420
$ws.Range("M65").Value = -11820
$ws.Range("N65").Value = -166685660

$ws.Range("H86").Value = 6689519
$ws.Range("I86").Value = 13336134
$ws.Range("J86").Value = 42904.2
$ws.Range("K86").Value = 13336134
$ws.Range("L86").Value = 42904.2
$ws.Range("M86").Value = -13335011
$ws.Range("N86").Value = -45150.2

$ws.Range("H89").Value = 6689519
$ws.Range("I89").Value = 13336134
$ws.Range("J89").Value = 42904.2
$ws.Range("K89").Value = 66680670
$ws.Range("L89").Value = 214521
$ws.Range("M89").Value = -66675054
$ws.Range("N89").Value = -225753

$ws.Range("H132").Value = 6703.48
$ws.Range("I132").Value = 8092.647
$ws.Range("J132").Value = 3751.5
$ws.Range("K132").Value = 24277.941
$ws.Range("L132").Value = 11254.5
$ws.Range("M132").Value = -21747.941
$ws.Range("N132").Value = -16314.5

$ws.Range("H135").Value = 49266.668
$ws.Range("J135").Value = 49266.668
$ws.Range("L135").Value = 49266.668
$ws.Range("N135").Value = -59406.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 792.6111
$ws.Range("I6").Value = 221.44444
$ws.Range("J6").Value = 1363.7778
$ws.Range("K6").Value = 664.33332
$ws.Range("L6").Value = 4091.3334
$ws.Range("M6").Value = -551.33332
$ws.Range("N6").Value = -4317.3334

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H68").Value = 2297.0527
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 2461.4119
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 7384.2357
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -9006.235700000001

$ws.Range("H71").Value = 2297.0527
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 2461.4119
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 22152.7071
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -30264.7071

$ws.Range("H131").Value = 13514720
$ws.Range("J131").Value = 1279.2616
$ws.Range("L131").Value = 3837.7848
$ws.Range("N131").Value = -13917.7848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1828500.8
$ws.Range("I12").Value = 1828500.8
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1828500.8
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1828360.8
$ws.Range("N12").ClearContents()

$ws.Range("H55").Value = 3312
$ws.Range("J55").Value = 3312
$ws.Range("L55").Value = 3312
$ws.Range("N55").Value = -3966

$ws.Range("H133").Value = 40644.5
$ws.Range("J133").Value = 40644.5
$ws.Range("L133").Value = 40644.5
$ws.Range("N133").Value = -50764.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 3000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3722

$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 3000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5496

$ws.Range("H108").Value = 25125.5
$ws.Range("J108").Value = 25125.5
$ws.Range("L108").Value = 25125.5
$ws.Range("N108").Value = -32805.5

$ws.Range("H132").Value = 36062.656
$ws.Range("I132").Value = 1275.7916
$ws.Range("J132").Value = 203039.6
$ws.Range("K132").Value = 3827.3748
$ws.Range("L132").Value = 609118.8
$ws.Range("M132").Value = -1297.3748
$ws.Range("N132").Value = -614178.8

$ws.Range("H135").Value = 31918.428
$ws.Range("J135").Value = 31918.428
$ws.Range("L135").Value = 31918.428
$ws.Range("N135").Value = -42058.428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7353.6665
$ws.Range("I41").Value = 5842
$ws.Range("J41").Value = 10377
$ws.Range("K41").Value = 5842
$ws.Range("L41").Value = 10377
$ws.Range("M41").Value = -5452
$ws.Range("N41").Value = -11157

$ws.Range("H45").Value = 8387.666999999999
$ws.Range("J45").Value = 8387.666999999999
$ws.Range("L45").Value = 8387.666999999999
$ws.Range("N45").Value = -9369.666999999999

$ws.Range("H122").Value = 14449140
$ws.Range("I122").Value = 15299048
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 45897144
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -45894694
$ws.Range("N122").Value = -7000

$ws.Range("H123").Value = 53991.855
$ws.Range("J123").Value = 53991.855
$ws.Range("L123").Value = 53991.855
$ws.Range("N123").Value = -63791.855

$ws.Range("H126").Value = 37038160
$ws.Range("I126").Value = 61729068
$ws.Range("J126").Value = 1794.5
$ws.Range("K126").Value = 185187204
$ws.Range("L126").Value = 5383.5
$ws.Range("M126").Value = -185184734
$ws.Range("N126").Value = -10323.5
